$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.29"
$ws.Range("E2").Value = "'1.48%"
$ws.Range("G2").Value = "'23"
$ws.Range("D3").Value = "'31.05"
$ws.Range("E3").Value = "'0.46%"
$ws.Range("G3").Value = "'23"
$ws.Range("D4").Value = "'4.933"
$ws.Range("E4").Value = "'1.32%"
$ws.Range("G4").Value = "'23"
$ws.Range("D5").Value = "'0.07370"
$ws.Range("E5").Value = "'2.92%"
$ws.Range("G5").Value = "'23"
$ws.Range("D6").Value = "'2.324"
$ws.Range("E6").Value = "'26.56%"
$ws.Range("G6").Value = "'23"
$ws.Range("D7").Value = "'7.692"
$ws.Range("E7").Value = "'0.59%"
$ws.Range("G7").Value = "'23"
$ws.Range("D8").Value = "'3.765"
$ws.Range("E8").Value = "'-0.26%"
$ws.Range("G8").Value = "'23"
$ws.Range("D9").Value = "'0.9123"
$ws.Range("E9").Value = "'2.01%"
$ws.Range("G9").Value = "'23"
$ws.Range("D10").Value = "'0.1689"
$ws.Range("E10").Value = "'2.67%"
$ws.Range("G10").Value = "'23"
$ws.Range("D11").Value = "'0.08187"
$ws.Range("G11").Value = "'23"
$ws.Range("D12").Value = "'0.08275"
$ws.Range("E12").Value = "'2.80%"
$ws.Range("G12").Value = "'23"
$ws.Range("D13").Value = "'0.03119"
$ws.Range("E13").Value = "'4.59%"
$ws.Range("G13").Value = "'23"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("G14").Value = "'23"
$ws.Range("D15").Value = "'0.001514"
$ws.Range("E15").Value = "'1.21%"
$ws.Range("G15").Value = "'23"
$ws.Range("D16").Value = "'0.005704"
$ws.Range("E16").Value = "'-2.04%"
$ws.Range("G16").Value = "'23"
$ws.Range("E17").Value = "'0.39%"
$ws.Range("G17").Value = "'23"
$ws.Range("D18").Value = "'2.075"
$ws.Range("E18").Value = "'-1.59%"
$ws.Range("G18").Value = "'23"
$ws.Range("G19").Value = "'23"
$ws.Range("E20").Value = "'0.88%"
$ws.Range("G20").Value = "'23"
$ws.Range("D21").Value = "'3.980"
$ws.Range("E21").Value = "'-6.67%"
$ws.Range("G21").Value = "'23"
$ws.Range("E22").Value = "'5.01%"
$ws.Range("G22").Value = "'23"
$ws.Range("D23").Value = "'0.04546"
$ws.Range("E23").Value = "'1.58%"
$ws.Range("G23").Value = "'23"
$ws.Range("E24").Value = "'-0.20%"
$ws.Range("G24").Value = "'23"
$ws.Range("D25").Value = "'0.004327"
$ws.Range("E25").Value = "'-7.30%"
$ws.Range("G25").Value = "'23"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("G26").Value = "'23"
$ws.Range("G27").Value = "'23"
$ws.Range("G28").Value = "'23"
$ws.Range("G29").Value = "'23"
$ws.Range("G30").Value = "'23"
$ws.Range("G31").Value = "'23"
$ws.Range("G32").Value = "'23"
$ws.Range("G33").Value = "'23"
$ws.Range("G34").Value = "'23"
$ws.Range("G35").Value = "'23"
$ws.Range("G36").Value = "'23"
$ws.Range("G37").Value = "'23"
$ws.Range("G38").Value = "'23"
$ws.Range("D39").Value = "'0.01606"
$ws.Range("G39").Value = "'23"
$ws.Range("D40").Value = "'0.04443"
$ws.Range("E40").Value = "'2.07%"
$ws.Range("G40").Value = "'23"
$ws.Range("D41").Value = "'0.007322"
$ws.Range("E41").Value = "'-0.79%"
$ws.Range("G41").Value = "'23"
$ws.Range("D42").Value = "'0.008847"
$ws.Range("G42").Value = "'23"
$ws.Range("D43").Value = "'0.1326"
$ws.Range("E43").Value = "'1.47%"
$ws.Range("G43").Value = "'23"
$ws.Range("D44").Value = "'0.002091"
$ws.Range("E44").Value = "'4.39%"
$ws.Range("G44").Value = "'23"
$ws.Range("D45").Value = "'0.009099"
$ws.Range("E45").Value = "'-11.30%"
$ws.Range("G45").Value = "'23"
$ws.Range("D46").Value = "'0.00005912"
$ws.Range("E46").Value = "'0.76%"
$ws.Range("G46").Value = "'23"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("G47").Value = "'23"
$ws.Range("E48").Value = "'1.27%"
$ws.Range("G48").Value = "'23"
$ws.Range("G49").Value = "'23"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("G50").Value = "'23"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("G51").Value = "'23"
